$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.419.48'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '3.139.57'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '610.43'
$ws.Range("E5").Value = '  +1.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.46'
$ws.Range("E6").Value = '  -1.86%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.136.33'
$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("E9").Value = '  +0.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  +0.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.35'
$ws.Range("E11").Value = '  -2.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.475'
$ws.Range("E12").Value = '  +0.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000257'
$ws.Range("E13").Value = '  +2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.54'
$ws.Range("E14").Value = '  -0.95%  '

$ws.Range("D15").Value = '3.654.14'
$ws.Range("E15").Value = '  -0.10%  '

$ws.Range("E16").Value = '  +2.78%  '

$ws.Range("D17").Value = '64.376.27'
$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("D18").Value = '3.138.80'
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.16'
$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.69'
$ws.Range("E21").Value = '  +0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.722'
$ws.Range("E22").Value = '  +2.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.81'
$ws.Range("E23").Value = '  +1.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.60'
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.00'
$ws.Range("E25").Value = '  +2.12%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.79'
$ws.Range("E27").Value = '  -2.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.54'
$ws.Range("E28").Value = '  +2.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.38'
$ws.Range("E29").Value = '  +9.89%  '

$ws.Range("E30").Value = '  +4.36%  '

$ws.Range("E31").Value = '  -4.33%  '

$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.69'
$ws.Range("E33").Value = '  +2.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.63'
$ws.Range("E34").Value = '  -3.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  +0.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.93'
$ws.Range("E36").Value = '  -0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.48'
$ws.Range("E37").Value = '  -3.12%  '

$ws.Range("D38").Value = '0.0₃0740'
$ws.Range("E38").Value = '  +4.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '452.66'
$ws.Range("E39").Value = '  +2.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.99'
$ws.Range("E40").Value = '  +4.05%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0396'
$ws.Range("E41").Value = '  +0.67%  '

$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.32'
$ws.Range("E43").Value = '  -0.95%  '

$ws.Range("D44").Value = '2.856.77'
$ws.Range("E44").Value = '  +1.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.265'
$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.26'
$ws.Range("E46").Value = '  +0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.42'
$ws.Range("E47").Value = '  +5.69%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.39'
$ws.Range("E48").Value = '  +0.47%  '

$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").Value = '  +0.14%  '

$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.85'
$ws.Range("E51").Value = '  +1.97%  '
